# Fruta / hortaliza, semanal
# Insert a new week's worth of data (3 rows: Especial/Primera/Segunda) right
# before the former row 894, pushing the rest of the table down by 3 rows
# (and the 3 rows that fall off the end land as new rows 947-949).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 894 downward (through the end of the used range) by 3 rows.
$ws.Rows("894:896").Insert()

# Populate the newly inserted rows with the new week's data.
$ws.Range("A894").Value = 3
$ws.Range("B894").Value = "Femacal de La Calera"
$ws.Range("C894").Value = "Coquimbo"
$ws.Range("D894").Value = 45021
$ws.Range("E894").Value = 5
$ws.Range("F894").Value = "Fruta"
$ws.Range("G894").Value = 100101
$ws.Range("H894").Value = "Berries"
$ws.Range("I894").Value = 100101007
$ws.Range("J894").Value = "Kiwi"
$ws.Range("K894").Value = "Hayward"
$ws.Range("L894").Value = "Especial"
$ws.Range("M894").Value = 75
$ws.Range("N894").Value = 10000
$ws.Range("O894").Value = 10000
$ws.Range("P894").Value = 10000
$ws.Range("Q894").Value = "`$/bandeja 10 kilos"
$ws.Range("R894").Value = "Región de O'Higgins"
$ws.Range("S894").Value = 1000
$ws.Range("T894").Value = 10

$ws.Range("A895").Value = 3
$ws.Range("B895").Value = "Femacal de La Calera"
$ws.Range("C895").Value = "Coquimbo"
$ws.Range("D895").Value = 45021
$ws.Range("E895").Value = 5
$ws.Range("F895").Value = "Fruta"
$ws.Range("G895").Value = 100101
$ws.Range("H895").Value = "Berries"
$ws.Range("I895").Value = 100101007
$ws.Range("J895").Value = "Kiwi"
$ws.Range("K895").Value = "Hayward"
$ws.Range("L895").Value = "Primera"
$ws.Range("M895").Value = 80
$ws.Range("N895").Value = 8000
$ws.Range("O895").Value = 8000
$ws.Range("P895").Value = 8000
$ws.Range("Q895").Value = "`$/bandeja 10 kilos"
$ws.Range("R895").Value = "Región de O'Higgins"
$ws.Range("S895").Value = 800
$ws.Range("T895").Value = 10

$ws.Range("A896").Value = 3
$ws.Range("B896").Value = "Femacal de La Calera"
$ws.Range("C896").Value = "Coquimbo"
$ws.Range("D896").Value = 45021
$ws.Range("E896").Value = 5
$ws.Range("F896").Value = "Fruta"
$ws.Range("G896").Value = 100101
$ws.Range("H896").Value = "Berries"
$ws.Range("I896").Value = 100101007
$ws.Range("J896").Value = "Kiwi"
$ws.Range("K896").Value = "Hayward"
$ws.Range("L896").Value = "Segunda"
$ws.Range("M896").Value = 70
$ws.Range("N896").Value = 7000
$ws.Range("O896").Value = 7000
$ws.Range("P896").Value = 7000
$ws.Range("Q896").Value = "`$/bandeja 10 kilos"
$ws.Range("R896").Value = "Región de O'Higgins"
$ws.Range("S896").Value = 700
$ws.Range("T896").Value = 10
